$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query text in B2 (Cases query): append an ORDER BY / LIMIT clause ---
$b2 = $ws.Range("B2").Text
$b2New = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Range("B2").Value = $b2New

# --- Update query text in B3 (Samples query): append an order By / LIMIT clause ---
$b3 = $ws.Range("B3").Text
$b3New = $b3 + "`norder By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $b3New

# --- Update query text in B4 (Files query): fix casing/spacing on ORDER BY clause and add LIMIT ---
$b4 = $ws.Range("B4").Text
$b4New = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value = $b4New

# --- Row heights grew (wrap-text autofit) now that B2/B3 hold more text ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# --- Update the active selection / scrolled view to cell B4 ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select() | Out-Null
